# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Terminal Hortofrutícola Agro Chillán - Pera"
# at rows 91-92, pushing the existing rows 91..123 down to 93..125.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 91 (shifts old rows 91-123 down to 93-125)
$ws.Rows("91:92").Insert()

# New row 91: Pera, Packham's Triumph, Especial, 60 units, week of 44455
$ws.Cells.Item(91, 1).Value = 7
$ws.Cells.Item(91, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(91, 3).Value = "Ñuble"
$ws.Cells.Item(91, 4).Value = 44455
$ws.Cells.Item(91, 5).Value = 16
$ws.Cells.Item(91, 6).Value = "Fruta"
$ws.Cells.Item(91, 7).Value = 100104
$ws.Cells.Item(91, 8).Value = "Frutos de pepita"
$ws.Cells.Item(91, 9).Value = 100104005
$ws.Cells.Item(91, 10).Value = "Pera"
$ws.Cells.Item(91, 11).Value = "Packham's Triumph"
$ws.Cells.Item(91, 12).Value = "Especial"
$ws.Cells.Item(91, 13).Value = 60
$ws.Cells.Item(91, 14).Value = 10000
$ws.Cells.Item(91, 15).Value = 11000
$ws.Cells.Item(91, 16).Value = 10500
$ws.Cells.Item(91, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(91, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(91, 19).Value = 656
$ws.Cells.Item(91, 20).Value = 16

# New row 92: Pera, Packham's Triumph, Primera, 60 units, week of 44455
$ws.Cells.Item(92, 1).Value = 7
$ws.Cells.Item(92, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(92, 3).Value = "Ñuble"
$ws.Cells.Item(92, 4).Value = 44455
$ws.Cells.Item(92, 5).Value = 16
$ws.Cells.Item(92, 6).Value = "Fruta"
$ws.Cells.Item(92, 7).Value = 100104
$ws.Cells.Item(92, 8).Value = "Frutos de pepita"
$ws.Cells.Item(92, 9).Value = 100104005
$ws.Cells.Item(92, 10).Value = "Pera"
$ws.Cells.Item(92, 11).Value = "Packham's Triumph"
$ws.Cells.Item(92, 12).Value = "Primera"
$ws.Cells.Item(92, 13).Value = 60
$ws.Cells.Item(92, 14).Value = 9200
$ws.Cells.Item(92, 15).Value = 9600
$ws.Cells.Item(92, 16).Value = 9400
$ws.Cells.Item(92, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(92, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(92, 19).Value = 588
$ws.Cells.Item(92, 20).Value = 16

# Re-apply the date number format (style s="2") to the new D91/D92 cells
$ws.Range("D91:D92").NumberFormat = "YYYY-MM-DD HH:MM:SS"
